$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Convert startdate/enddate (C2:D5) from real dates to literal text ---
$ws.Range("C2:D5").NumberFormat = "@"

$ws.Range("C2").Value = "2021-09-07"
$ws.Range("C3").Value = "2021-09-07"
$ws.Range("C4").Value = "2021-09-07"
$ws.Range("C5").Value = "2021-09-07"

$ws.Range("D2").Value = "2022-06-07"
$ws.Range("D3").Value = "2022-06-07"
$ws.Range("D4").Value = "2022-06-07"
$ws.Range("D5").Value = "2022-06-07"

# --- Give the header row (A1:D1) the same text-number-format style used
#     by the date columns, by copying the cell format (not the value). ---
$ws.Range("C2").Copy()
$ws.Range("A1:D1").PasteSpecial(-4122)  # xlPasteFormats
$null = $excel.CutCopyMode

# --- Move the saved cursor position to C9 ---
$null = $ws.Range("C9").Select()

# --- Touch Page Setup so an explicit <pageSetup orientation="portrait".../>
#     gets persisted on save. ---
$ws.PageSetup.Orientation = 1
